$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7143781781196594
$ws.Range("B1").Value = 0.8832109570503235
$ws.Range("C1").Value = 0.7768504619598389
$ws.Range("D1").Value = 3.347337484359741
$ws.Range("E1").Value = 1.635702967643738
